# Shopping slots export: drop "대행사" (agency) and "상태" (status) columns,
# and seed two sample rows, per the new agency-managed slot quota feature.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "대행사" column (old column B). Everything shifts left.
$ws.Range("B1").EntireColumn.Delete()

# Remove the "상태" column. After the deletion above it now sits at column K.
$ws.Range("K1").EntireColumn.Delete()

# --- Sample row 1 ---
$ws.Range("A2").Value = "샘플 슬롯 1"
$ws.Range("B2").Value = "스마트스토어"
$ws.Range("C2").Value = "P12345"
$ws.Range("D2").Value = "샘플 상품 1"
$ws.Range("E2").Value = "키워드1,키워드2"
$ws.Range("F2").Value = 10000
$ws.Range("G2").Value = 8000
$ws.Range("H2").Value = 45778
$ws.Range("H2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("H2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("I2").Value = 45808
$ws.Range("I2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J2").Value = "CPC"
$ws.Range("K2").Value = 5000

# --- Sample row 2 ---
$ws.Range("A3").Value = "샘플 슬롯 2"
$ws.Range("B3").Value = "브랜드몰"
$ws.Range("C3").Value = "P67890"
$ws.Range("D3").Value = "샘플 상품 2"
$ws.Range("E3").Value = "키워드3,키워드4"
$ws.Range("F3").Value = 20000
$ws.Range("G3").Value = 18000
$ws.Range("H3").Value = 45792
$ws.Range("H3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("I3").Value = 45823
$ws.Range("I3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J3").Value = "CPM"
$ws.Range("K3").Value = 6000
